$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the oldest quarter (column D, "فصل دوم منتهی به 1399/06") -- this
# shifts every later quarter's column one position to the left (E->D, ... M->L)
$ws.Columns("D").Delete()

# Bring in a new column M for the newest quarter by cloning the formatting
# (styles/column width) of the now-last data column L
$ws.Range("L1:L28").Copy()
$ws.Range("M1:M28").PasteSpecial(-4122)
$ws.Columns("M").ColumnWidth = 30.166666666666668

# Headers for the new quarter column
$ws.Range("M8").Value = "فصل چهارم منتهی به 1401/12"

# "1402-02-28" looks like a date to Excel's auto-detection, so force text
# and then re-pull the formatting (border/fill/font) back from the sibling
# date cell L9, which the NumberFormat change above disturbs
$ws.Range("M9").NumberFormat = "@"
$ws.Range("M9").Value = "1402-02-28"
$ws.Range("L9").Copy()
$ws.Range("M9").PasteSpecial(-4122)

# The quarter that used to sit in column J ("فصل چهارم منتهی به 1400/12")
# got a later amendment -- update its publish-date label
$ws.Range("I9").Value = "1402-02-28 (8)"

# New financial figures for the newest quarter (column M)
$ws.Range("M11").Value = 4565621
$ws.Range("M12").Value = -2631918
$ws.Range("M13").Value = 1933703
$ws.Range("M14").Value = -178886
$ws.Range("M15").Value = 0
$ws.Range("M16").Value = 270851
$ws.Range("M17").Value = 2025668
$ws.Range("M18").Value = -139315
$ws.Range("M19").Value = -13850
$ws.Range("M20").Value = 1872503
$ws.Range("M21").Value = -158951
$ws.Range("M22").Value = 1713552
$ws.Range("M23").Value = 0
$ws.Range("M24").Value = 1713552
$ws.Range("M25").Value = 2397
$ws.Range("M26").Value = 715000
$ws.Range("M27").Value = 2397
